$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 628, shifting existing rows 628-714 down to 630-716.
$ws.Rows("628:629").Insert()

# New row 628 values
$ws.Range("A628").Value = 10
$ws.Range("B628").Value = "Vega Modelo de Temuco"
$ws.Range("C628").Value = "La Araucanía"
$ws.Range("D628").Value = 45124
$ws.Range("E628").Value = 9
$ws.Range("F628").Value = 100112023
$ws.Range("G628").Value = "Brócoli"
$ws.Range("H628").Value = "Sin especificar"
$ws.Range("I628").Value = "Primera"
$ws.Range("J628").Value = 3000
$ws.Range("K628").Value = 1000
$ws.Range("L628").Value = 1000
$ws.Range("M628").Value = 1000
$ws.Range("N628").Value = "$/unidad"
$ws.Range("O628").Value = "Región Metropolitana"
$ws.Range("P628").Value = 1000
$ws.Range("Q628").Value = 1
$ws.Range("R628").Value = "Hortaliza"

# New row 629 values
$ws.Range("A629").Value = 10
$ws.Range("B629").Value = "Vega Modelo de Temuco"
$ws.Range("C629").Value = "La Araucanía"
$ws.Range("D629").Value = 45124
$ws.Range("E629").Value = 9
$ws.Range("F629").Value = 100112023
$ws.Range("G629").Value = "Brócoli"
$ws.Range("H629").Value = "Sin especificar"
$ws.Range("I629").Value = "Primera"
$ws.Range("J629").Value = 2000
$ws.Range("K629").Value = 1300
$ws.Range("L629").Value = 1300
$ws.Range("M629").Value = 1300
$ws.Range("N629").Value = "$/unidad"
$ws.Range("O629").Value = "Región del Maule"
$ws.Range("P629").Value = 1300
$ws.Range("Q629").Value = 1
$ws.Range("R629").Value = "Hortaliza"

# Ensure the date cells keep the date/time number format used by the rest of column D.
$ws.Range("D628:D629").NumberFormat = $ws.Range("D630").NumberFormat
